# Ore aggiornato al 08/01 al pomeriggio
#
# Adds two entries for 08/01/2017 (afternoon) - one for Giovanni, one for
# Mirko - describing the discussion about implementations, making methods
# virtual, adding missing signatures and personal branches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$activityDate = 42743  # 08/01/2017
$hours = 0.14583333333333334  # 3.5 hours
$desc = "discussione implementazioni e sistemato virtual in tutti I metodi. Aggiunte firme alle classi mancanti. Aggiunti branch personali."

# New row 22 - Giovanni
$ws.Cells.Item(21, 1).Copy()
$ws.Cells.Item(22, 1).PasteSpecial(-4122)
$ws.Cells.Item(21, 3).Copy()
$ws.Cells.Item(22, 3).PasteSpecial(-4122)
$ws.Cells.Item(21, 4).Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4122)

$ws.Cells.Item(22, 1).Value = $activityDate
$ws.Cells.Item(22, 2).Value = "Giovanni"
$ws.Cells.Item(22, 3).Value = $desc
$ws.Cells.Item(22, 4).Value = $hours

# New row 23 - Mirko
$ws.Cells.Item(21, 1).Copy()
$ws.Cells.Item(23, 1).PasteSpecial(-4122)
$ws.Cells.Item(21, 3).Copy()
$ws.Cells.Item(23, 3).PasteSpecial(-4122)
$ws.Cells.Item(21, 4).Copy()
$ws.Cells.Item(23, 4).PasteSpecial(-4122)

$ws.Cells.Item(23, 1).Value = $activityDate
$ws.Cells.Item(23, 2).Value = "Mirko"
$ws.Cells.Item(23, 3).Value = $desc
$ws.Cells.Item(23, 4).Value = $hours

$excel.CutCopyMode = $false

# The long wrapped description needs a taller row to display fully.
$ws.Rows.Item(22).RowHeight = 87.45
$ws.Rows.Item(23).RowHeight = 87.45

# Scroll back to the top of the sheet and leave the grand-total cell selected.
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I2").Select() | Out-Null

$wb.Save() | Out-Null
